# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 151
$wsExhibit.Range("F4").Value = 241
$wsExhibit.Range("F5").Value = 3863
$wsExhibit.Range("F6").Value = 27
$wsExhibit.Range("F7").Value = 436

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 151
$wsAll.Range("F4").Value = 241
$wsAll.Range("F5").Value = 3863
$wsAll.Range("F8").Value = 27
$wsAll.Range("F9").Value = 436
